$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B/C/E columns (and non-numeric-looking D values) are safe to assign directly;
# D-column numeric-looking strings need a text-prefix so Excel keeps them as text
# (matching the original inlineStr/Text cell type) instead of coercing to a Number.

$ws.Range('D2').Value = '35.578.19'
$ws.Range('E2').Value = '  -2.76%  '
$ws.Range('D3').Value = '1.983.73'
$ws.Range('E3').Value = '  -3.67%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'245.17"
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('D6').Value = "'0.637"
$ws.Range('E6').Value = '  -4.50%  '
$ws.Range('D7').Value = "'57.69"
$ws.Range('E7').Value = '  +5.92%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = "'58.40"
$ws.Range('E9').Value = '  -0.38%  '
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').Value = "'0.0735"
$ws.Range('E11').Value = '  -2.17%  '
$ws.Range('E12').Value = '  -2.77%  '
$ws.Range('D13').Value = "'0.945"
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').Value = "'14.45"
$ws.Range('E14').Value = '  -1.86%  '
$ws.Range('D15').Value = '2.270.12'
$ws.Range('E15').Value = '  -3.85%  '
$ws.Range('D16').Value = "'5.29"
$ws.Range('E16').Value = '  -2.60%  '
$ws.Range('D17').Value = '1.978.34'
$ws.Range('E17').Value = '  -4.94%  '
$ws.Range('D18').Value = "'17.86"
$ws.Range('E18').Value = '  +6.46%  '
$ws.Range('D19').Value = '35.583.24'
$ws.Range('E19').Value = '  -2.57%  '
$ws.Range('E20').Value = '  -1.04%  '
$ws.Range('E21').Value = '  -1.81%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'5.17"
$ws.Range('E22').Value = '  -1.76%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').Value = "'232.96"
$ws.Range('E23').Value = '  -2.14%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('E25').Value = '  +19.83%  '
$ws.Range('E26').Value = '  -2.99%  '
$ws.Range('D27').Value = "'164.75"
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').Value = "'9.15"
$ws.Range('E28').Value = '  -1.70%  '
$ws.Range('D29').Value = "'19.20"
$ws.Range('E29').Value = '  -4.60%  '
$ws.Range('E30').Value = '  -2.54%  '
$ws.Range('E31').Value = '  -4.04%  '
$ws.Range('E32').Value = '  -7.42%  '
$ws.Range('E33').Value = '  +16.40%  '
$ws.Range('E34').Value = '  -0.80%  '
$ws.Range('E35').Value = '  +9.73%  '
$ws.Range('D36').Value = "'4.33"
$ws.Range('E36').Value = '  -3.71%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E38').Value = '  -3.79%  '
$ws.Range('D39').Value = "'5.26"
$ws.Range('E39').Value = '  +8.75%  '
$ws.Range('E40').Value = '  -2.51%  '
$ws.Range('E41').Value = '  +1.12%  '
$ws.Range('E42').Value = '  -2.13%  '
$ws.Range('E43').Value = '  -1.74%  '
$ws.Range('D44').Value = "'7.68"
$ws.Range('E44').Value = '  +0.83%  '
$ws.Range('D45').Value = "'91.95"
$ws.Range('E45').Value = '  -2.57%  '
$ws.Range('D46').Value = "'16.08"
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('D47').Value = "'0.0889"
$ws.Range('E47').Value = '  -2.54%  '
$ws.Range('D48').Value = '1.369.61'
$ws.Range('E48').Value = '  -3.46%  '
$ws.Range('E49').Value = '  +1.10%  '
$ws.Range('D50').Value = "'46.78"
$ws.Range('E50').Value = '  +3.37%  '
$ws.Range('E51').Value = '  -1.70%  '
